$d = $word.ActiveDocument

function FindParagraphIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -like "*$text*") {
            return $i
        }
    }
    return -1
}

$idx = FindParagraphIndex("No good, JRE only 1.5 even in later versions")

# Insert first new bullet right after the "No good..." paragraph.
$anchor = $d.Paragraphs($idx)
$anchor.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs($idx + 1)
$p1.Range.InsertBefore("Do analysis ourselves")
$p1.Range.ListFormat.ListLevelNumber = 3

# Insert second new bullet right after the first one.
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($idx + 2)
$p2.Range.InsertBefore("Problem: type verification will trigger class loading, which needs a current thread, which may not exist in control flow of getMirageClass()")
$p2.Range.ListFormat.ListLevelNumber = 3
